# Auto-generated script applying the scheduled Golem Profits price/profit refresh
# across the ALC/ARM/BSM/CUL/GSM/LTW/WVR sheets (each an Excel Table of leve data).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 287.5
$ws.Cells.Item(2, 9).Value = 287.5
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 287.5
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = -174.5
$ws.Cells.Item(2, 14).ClearContents()  # N2: was -508
$ws.Cells.Item(41, 8).Value = 544.75
$ws.Cells.Item(41, 10).Value = 250
$ws.Cells.Item(41, 12).Value = 250
$ws.Cells.Item(41, 14).Value = -1130
$ws.Cells.Item(58, 8).Value = 3387.25
$ws.Cells.Item(58, 9).Value = 50
$ws.Cells.Item(58, 10).Value = 4499.6665
$ws.Cells.Item(58, 11).Value = 150
$ws.Cells.Item(58, 12).Value = 13498.9995
$ws.Cells.Item(58, 13).Value = 0
$ws.Cells.Item(58, 14).Value = -13798.9995
$ws.Cells.Item(64, 8).Value = 2638
$ws.Cells.Item(64, 10).Value = 2297
$ws.Cells.Item(64, 12).Value = 2297
$ws.Cells.Item(64, 14).Value = -2793
$ws.Cells.Item(67, 8).Value = 2638
$ws.Cells.Item(67, 10).Value = 2297
$ws.Cells.Item(67, 12).Value = 2297
$ws.Cells.Item(67, 14).Value = -4013
$ws.Cells.Item(76, 8).Value = 3899.5
$ws.Cells.Item(76, 9).Value = 3800
$ws.Cells.Item(76, 10).Value = 3999
$ws.Cells.Item(76, 11).Value = 3800
$ws.Cells.Item(76, 12).Value = 3999
$ws.Cells.Item(76, 13).Value = -3485
$ws.Cells.Item(76, 14).Value = -4629
$ws.Cells.Item(79, 8).Value = 3899.5
$ws.Cells.Item(79, 9).Value = 3800
$ws.Cells.Item(79, 10).Value = 3999
$ws.Cells.Item(79, 11).Value = 3800
$ws.Cells.Item(79, 12).Value = 3999
$ws.Cells.Item(79, 13).Value = -2708
$ws.Cells.Item(79, 14).Value = -6183
$ws.Cells.Item(105, 8).Value = 63333.332
$ws.Cells.Item(105, 10).Value = 63333.332
$ws.Cells.Item(105, 12).Value = 63333.332
$ws.Cells.Item(105, 14).Value = -70321.33199999999
$ws.Cells.Item(115, 8).Value = 1278.75
$ws.Cells.Item(115, 9).Value = 1278.75
$ws.Cells.Item(115, 11).Value = 3836.25
$ws.Cells.Item(115, 13).Value = -2269.25
$ws.Cells.Item(131, 8).Value = 7777
$ws.Cells.Item(131, 9).Value = 7777
$ws.Cells.Item(131, 11).Value = 23331
$ws.Cells.Item(131, 13).Value = -18291

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(25, 8).Value = 2050
$ws.Cells.Item(25, 9).Value = 2050
$ws.Cells.Item(25, 11).Value = 2050
$ws.Cells.Item(25, 13).Value = -1648
$ws.Cells.Item(32, 8).Value = 817.3333
$ws.Cells.Item(32, 9).Value = 469.625
$ws.Cells.Item(32, 11).Value = 469.625
$ws.Cells.Item(32, 13).Value = -182.625
$ws.Cells.Item(35, 8).Value = 34995
$ws.Cells.Item(35, 9).Value = 0
$ws.Cells.Item(35, 10).Value = 34995
$ws.Cells.Item(35, 11).Value = 0
$ws.Cells.Item(35, 12).Value = 34995
$ws.Cells.Item(35, 13).ClearContents()  # M35: was -393
$ws.Cells.Item(35, 14).Value = -35807
$ws.Cells.Item(63, 8).Value = 32500
$ws.Cells.Item(63, 9).Value = 15000
$ws.Cells.Item(63, 11).Value = 15000
$ws.Cells.Item(63, 13).Value = -14314
$ws.Cells.Item(66, 8).Value = 32500
$ws.Cells.Item(66, 9).Value = 15000
$ws.Cells.Item(66, 11).Value = 75000
$ws.Cells.Item(66, 13).Value = -71568
$ws.Cells.Item(88, 8).Value = 1039.7273
$ws.Cells.Item(88, 9).Value = 598.5714
$ws.Cells.Item(88, 10).Value = 1811.75
$ws.Cells.Item(88, 11).Value = 598.5714
$ws.Cells.Item(88, 12).Value = 1811.75
$ws.Cells.Item(88, 13).Value = -192.5714
$ws.Cells.Item(88, 14).Value = -2623.75
$ws.Cells.Item(91, 8).Value = 1039.7273
$ws.Cells.Item(91, 9).Value = 598.5714
$ws.Cells.Item(91, 10).Value = 1811.75
$ws.Cells.Item(91, 11).Value = 598.5714
$ws.Cells.Item(91, 12).Value = 1811.75
$ws.Cells.Item(91, 13).Value = 805.4286
$ws.Cells.Item(91, 14).Value = -4619.75
$ws.Cells.Item(102, 8).Value = 1388.4546
$ws.Cells.Item(102, 9).Value = 626.625
$ws.Cells.Item(102, 10).Value = 3420
$ws.Cells.Item(102, 11).Value = 626.625
$ws.Cells.Item(102, 12).Value = 3420
$ws.Cells.Item(102, 13).Value = 995.375
$ws.Cells.Item(102, 14).Value = -6664

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 459.14285
$ws.Cells.Item(80, 9).Value = 219.14285
$ws.Cells.Item(80, 10).Value = 699.1429000000001
$ws.Cells.Item(80, 11).Value = 219.14285
$ws.Cells.Item(80, 12).Value = 699.1429000000001
$ws.Cells.Item(80, 13).Value = 778.85715
$ws.Cells.Item(80, 14).Value = -2695.1429
$ws.Cells.Item(83, 8).Value = 459.14285
$ws.Cells.Item(83, 9).Value = 219.14285
$ws.Cells.Item(83, 10).Value = 699.1429000000001
$ws.Cells.Item(83, 11).Value = 1095.71425
$ws.Cells.Item(83, 12).Value = 3495.7145
$ws.Cells.Item(83, 13).Value = 3896.28575
$ws.Cells.Item(83, 14).Value = -13479.7145
$ws.Cells.Item(86, 8).Value = 1772.5
$ws.Cells.Item(86, 9).Value = 1696.6666
$ws.Cells.Item(86, 11).Value = 1696.6666
$ws.Cells.Item(86, 13).Value = -573.6666
$ws.Cells.Item(89, 8).Value = 1772.5
$ws.Cells.Item(89, 9).Value = 1696.6666
$ws.Cells.Item(89, 11).Value = 8483.333000000001
$ws.Cells.Item(89, 13).Value = -2867.333000000001
$ws.Cells.Item(107, 8).Value = 1977.6
$ws.Cells.Item(107, 9).Value = 1977.6
$ws.Cells.Item(107, 11).Value = 1977.6
$ws.Cells.Item(107, 13).Value = -57.59999999999991
$ws.Cells.Item(140, 8).Value = 130000
$ws.Cells.Item(140, 10).Value = 130000
$ws.Cells.Item(140, 12).Value = 130000
$ws.Cells.Item(140, 14).Value = -140360

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 349.2857
$ws.Cells.Item(2, 9).Value = 101
$ws.Cells.Item(2, 10).Value = 390.66666
$ws.Cells.Item(2, 11).Value = 606
$ws.Cells.Item(2, 12).Value = 2343.99996
$ws.Cells.Item(2, 13).Value = -493
$ws.Cells.Item(2, 14).Value = -2569.99996
$ws.Cells.Item(38, 8).Value = 483.47827
$ws.Cells.Item(38, 9).Value = 232.2
$ws.Cells.Item(38, 10).Value = 553.2778
$ws.Cells.Item(38, 11).Value = 696.5999999999999
$ws.Cells.Item(38, 12).Value = 1659.8334
$ws.Cells.Item(38, 13).Value = -349.5999999999999
$ws.Cells.Item(38, 14).Value = -2353.8334
$ws.Cells.Item(137, 8).Value = 7500
$ws.Cells.Item(137, 10).Value = 7500
$ws.Cells.Item(137, 12).Value = 22500
$ws.Cells.Item(137, 14).Value = -32700

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(104, 8).Value = 6500
$ws.Cells.Item(104, 10).Value = 6500
$ws.Cells.Item(104, 12).Value = 6500
$ws.Cells.Item(104, 14).Value = -13488
$ws.Cells.Item(113, 8).Value = 2966.3333
$ws.Cells.Item(113, 9).Value = 1949.75
$ws.Cells.Item(113, 10).Value = 4999.5
$ws.Cells.Item(113, 11).Value = 1949.75
$ws.Cells.Item(113, 12).Value = 4999.5
$ws.Cells.Item(113, 13).Value = 220.25
$ws.Cells.Item(113, 14).Value = -9339.5
$ws.Cells.Item(122, 8).Value = 2638.3076
$ws.Cells.Item(122, 9).Value = 2229.9
$ws.Cells.Item(122, 10).Value = 3999.6667
$ws.Cells.Item(122, 11).Value = 6689.700000000001
$ws.Cells.Item(122, 12).Value = 11999.0001
$ws.Cells.Item(122, 13).Value = -4239.700000000001
$ws.Cells.Item(122, 14).Value = -16899.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 558927.7
$ws.Cells.Item(40, 9).Value = 1891.6666
$ws.Cells.Item(40, 10).Value = 1672999.6
$ws.Cells.Item(40, 11).Value = 1891.6666
$ws.Cells.Item(40, 12).Value = 1672999.6
$ws.Cells.Item(40, 13).Value = -1755.6666
$ws.Cells.Item(40, 14).Value = -1673271.6
$ws.Cells.Item(46, 8).Value = 227386
$ws.Cells.Item(46, 9).Value = 336745.66
$ws.Cells.Item(46, 10).Value = 8666.666999999999
$ws.Cells.Item(46, 11).Value = 336745.66
$ws.Cells.Item(46, 12).Value = 8666.666999999999
$ws.Cells.Item(46, 13).Value = -336557.66
$ws.Cells.Item(46, 14).Value = -9042.666999999999
$ws.Cells.Item(61, 8).Value = 2392.7144
$ws.Cells.Item(61, 9).Value = 2392.7144
$ws.Cells.Item(61, 11).Value = 2392.7144
$ws.Cells.Item(61, 13).Value = -2190.7144
$ws.Cells.Item(82, 8).Value = 1490.5454
$ws.Cells.Item(82, 9).Value = 1210.4
$ws.Cells.Item(82, 11).Value = 1210.4
$ws.Cells.Item(82, 13).Value = -849.4000000000001
$ws.Cells.Item(85, 8).Value = 1490.5454
$ws.Cells.Item(85, 9).Value = 1210.4
$ws.Cells.Item(85, 11).Value = 1210.4
$ws.Cells.Item(85, 13).Value = 37.59999999999991
$ws.Cells.Item(113, 8).Value = 2392.7144
$ws.Cells.Item(113, 9).Value = 2392.7144
$ws.Cells.Item(113, 11).Value = 2392.7144
$ws.Cells.Item(113, 13).Value = -222.7143999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1806
$ws.Cells.Item(122, 9).Value = 1711
$ws.Cells.Item(122, 10).Value = 1996
$ws.Cells.Item(122, 11).Value = 5133
$ws.Cells.Item(122, 12).Value = 5988
$ws.Cells.Item(122, 13).Value = -2683
$ws.Cells.Item(122, 14).Value = -10888
